# act tablas web jul25
# Refresh "Gasto publico anual per capita en vivienda" data table (2023 update)
# and add an "actualizacion" metadata row.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsMeta = $wb.Worksheets.Item("Metadata")

# ---------------------------------------------------------------------------
# 1. "Data" sheet: refresh the Fecha/Valor series with the Jul-2025 dataset.
#    Two more recent years (2023, 2022) are prepended and every value is
#    re-estimated, so the whole A2:B35 block is rewritten.
# ---------------------------------------------------------------------------

$years = @(2023,2022,2021,2020,2019,2018,2017,2016,2015,2014,2013,2012,2011,2010,2009,2008,2007,2006,2005,2004,2003,2002,2001,2000,1999,1998,1997,1996,1995,1994,1993,1992,1991,1990)
$values = @(95.6,111.3,90.3,86.5,100.4,110.3,94.2,94.4,81.3,83.3,79.5,77.1,77.9,67.3,47.3,43.9,51.1,41.9,27.7,32.3,30.6,36.9,46.3,44.1,50.1,76,73.1,73.7,70.6,15.6,26.7,25.5,19.8,19.3)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $cellA = $wsData.Cells.Item($row, 1)
    # Force text storage so the year labels ("2023", "2022", ...) don't get
    # auto-coerced into numeric cells.
    $cellA.NumberFormat = "@"
    $cellA.Value = [string]$years[$i]
    $wsData.Cells.Item($row, 2).Value2 = $values[$i]
}

# ---------------------------------------------------------------------------
# 2. "Metadata" sheet: blank leading cell becomes a single space, and a new
#    "actualizacion" / "Julio 2025" row is inserted right before "cita".
# ---------------------------------------------------------------------------

$wsMeta.Range("A1").Value = " "

# "cita" currently sits on row 9 - push it (and everything after) down one
# row so we can drop the new pair in cleanly.
$wsMeta.Rows.Item(9).Insert()

$wsMeta.Cells.Item(9, 1).Value = "actualizacion"
$wsMeta.Cells.Item(9, 2).Value = "Julio 2025"
